$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find failed for: $old"
    }
}

# 1) Merge "Identifique pelo menos " + "3" + " ex  possíveis..." into a single run
#    (also removes the surrounding proofErr gramStart/gramEnd markers)
Replace-Exact "Identifique pelo menos 3 ex  possíveis de dados que você gerou hoje durante o dia." "Identifique pelo menos 3 ex  possíveis de dados que você gerou hoje durante o dia."

# 2) Append a new run with "." after "Já pensou em como aplicar ... Se sim, como?" (first occurrence,
#    inside the "Também queremos te conhecer" list). Force a transient property change so the engine
#    keeps it as a distinct <w:r> instead of silently re-merging it into the previous run.
$rng = $d.Content
$rng.Find.Execute("Já pensou em como aplicar Ciência de Dados para alcançar seus objetivos? Se sim, como?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.Text = "."
$rng.Bold = 1
$rng.Bold = 0

# 3) Merge " informações brutas..." + "  " + "sem contexto claro, " into a single run
Replace-Exact " informações brutas que, por si só, não têm muito valor ou significado, eles podem ser números, palavras ou observações  sem contexto claro, " " informações brutas que, por si só, não têm muito valor ou significado, eles podem ser números, palavras ou observações  sem contexto claro, "

# 4) Merge "Os dados são importantes na sociedade" + "  " + "contemporânea porque..." into one run
Replace-Exact "Os dados são importantes na sociedade  contemporânea porque influenciam decisões em diversas áreas, como saúde, educação, marketing e transporte, especialmente  após as pandemia da Covid-19, que destacou ainda mais sua relevância " "Os dados são importantes na sociedade  contemporânea porque influenciam decisões em diversas áreas, como saúde, educação, marketing e transporte, especialmente  após as pandemia da Covid-19, que destacou ainda mais sua relevância "

# 5) Merge "Estruturados " + "( " + "Tabelas)" into one run
Replace-Exact "Estruturados ( Tabelas)" "Estruturados ( Tabelas)"

# 6) Merge "Não estruturados (imagens, textos" + ")" into one run
Replace-Exact "Não estruturados (imagens, textos)" "Não estruturados (imagens, textos)"

# 7) "Semi estrutura" + "dos" -> "Semi estrutura" + "dos (O meio termo Inteligente, organizado)",
#    then reorder/merge the remaining runs of that sentence.
Replace-Exact "dos (O meio termo Inteligente, organizado)" "dos (O meio termo Inteligente, organizado)"
Replace-Exact "que não estão em tabelas , mas possuem marcadores/tags que dão alguma marcação" "que não estão em tabelas , mas possuem marcadores/tags que dão alguma marcação"

# 8) Merge "✅" + "  " + "2 - " into "  2 - " (single run), dropping proofErr markers
Replace-Exact "  2 - " "  2 - "

# 9) Merge "✅" + "  " + "3 - " into "  3 - " (single run), dropping proofErr markers
Replace-Exact "  3 - " "  3 - "

# 10) Merge "Segue também minha planilha no Excel com " + "3" + " abas..." into one run
Replace-Exact "Segue também minha planilha no Excel com 3 abas que criei para acompanhar o meu progresso em meu aprendizado no curso. Sendo a primeira de ‘Registro diário’, segunda do ‘Resumo’ para calcular a métrica e a terceira aba sendo’ Progresso Geral’." "Segue também minha planilha no Excel com 3 abas que criei para acompanhar o meu progresso em meu aprendizado no curso. Sendo a primeira de ‘Registro diário’, segunda do ‘Resumo’ para calcular a métrica e a terceira aba sendo’ Progresso Geral’."
